# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Wed Apr 17 19:48:40 UTC 2024 with GitHub Actions".
#
# Cells in column D hold plain text (e.g. "533.84", "60.928.17") even though
# many look numeric. Excel auto-converts a numeric-looking Value assignment
# into a real number (losing formatting like "1.00" -> 1), so any D-column
# value that parses as a number is written with a leading apostrophe to force
# text, exactly like typing '533.84 into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.928.17'
$ws.Range("E2").Value = '  -2.88%  '

$ws.Range("D3").Value = '2.993.32'
$ws.Range("E3").Value = '  -2.19%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''533.84'
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("D6").Value = '''133.53'
$ws.Range("E6").Value = '  +1.12%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '2.986.62'
$ws.Range("E8").Value = '  -2.20%  '

$ws.Range("E9").Value = '  +1.11%  '

$ws.Range("E10").Value = '  -3.44%  '

$ws.Range("D11").Value = '''6.09'
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").Value = '''0.444'
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("E13").Value = '  -1.12%  '

$ws.Range("D14").Value = '''34.06'
$ws.Range("E14").Value = '  +0.39%  '

$ws.Range("D15").Value = '3.481.56'
$ws.Range("E15").Value = '  -1.99%  '

$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").Value = '60.985.38'
$ws.Range("E17").Value = '  -2.87%  '

$ws.Range("D18").Value = '3.004.38'
$ws.Range("E18").Value = '  -1.91%  '

$ws.Range("D19").Value = '''6.59'
$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("D20").Value = '''462.47'
$ws.Range("E20").Value = '  -3.47%  '

$ws.Range("D21").Value = '''13.18'
$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("D22").Value = '''0.673'
$ws.Range("E22").Value = '  -2.20%  '

$ws.Range("E23").Value = '  -1.54%  '

$ws.Range("D24").Value = '''79.15'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("D25").Value = '''12.01'
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  -0.35%  '

$ws.Range("D28").Value = '''7.84'
$ws.Range("E28").Value = '  -2.05%  '

$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.21%  '

$ws.Range("D30").Value = '''1.89'
$ws.Range("E30").Value = '  +1.79%  '

$ws.Range("D31").Value = '''25.44'
$ws.Range("E31").Value = '  -1.45%  '

$ws.Range("D32").Value = '''1.14'
$ws.Range("E32").Value = '  +3.07%  '

$ws.Range("E33").Value = '  +2.88%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '''55.34'
$ws.Range("E34").Value = '  -2.65%  '

$ws.Range("B35").Value = 'Stacks'
$ws.Range("C35").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D35").Value = '''2.27'
$ws.Range("E35").Value = '  -2.69%  '

$ws.Range("D36").Value = '''5.87'
$ws.Range("E36").Value = '  -1.46%  '

$ws.Range("D37").Value = '''454.91'
$ws.Range("E37").Value = '  -3.87%  '

$ws.Range("D38").Value = '3.200.71'
$ws.Range("E38").Value = '  +3.80%  '

$ws.Range("D39").Value = '''0.0786'
$ws.Range("E39").Value = '  -0.32%  '

$ws.Range("D40").Value = '''0.0383'
$ws.Range("E40").Value = '  -1.96%  '

$ws.Range("E41").Value = '  +2.52%  '

$ws.Range("D42").Value = '''8.14'
$ws.Range("E42").Value = '  +1.44%  '

$ws.Range("D43").Value = '''27.52'
$ws.Range("E43").Value = '  +14.01%  '

$ws.Range("D44").Value = '''2.45'
$ws.Range("E44").Value = '  -5.27%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").Value = '''0.245'
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("D48").Value = '''118.98'
$ws.Range("E48").Value = '  -1.60%  '

$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("E50").Value = '  -7.18%  '

$ws.Range("D51").Value = '''1.25'
$ws.Range("E51").Value = '  +8.13%  '
